$wb = $excel.ActiveWorkbook

# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 534
    $ws.Range("F3").Value = 3487
    $ws.Range("F5").Value = 680
}
